$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-16 Monday" "2026-02-17 Tuesday"

Replace-Text "434×4=1736" "357×6=2142"
Replace-Text "183×3=549" "689×2=1378"
Replace-Text "973×2=1946" "810×5=4050"
Replace-Text "198×2=396" "407×5=2035"
Replace-Text "887×8=7096" "987×8=7896"

Replace-Text "417×7=2919" "623×2=1246"
Replace-Text "165×6=990" "590×9=5310"
Replace-Text "475×8=3800" "549×6=3294"
Replace-Text "297×9=2673" "752×2=1504"
Replace-Text "133×8=1064" "886×3=2658"

Replace-Text "866×3=2598" "609×6=3654"
Replace-Text "969×4=3876" "314×7=2198"
Replace-Text "271×4=1084" "470×8=3760"
Replace-Text "709×6=4254" "250×7=1750"
Replace-Text "973×9=8757" "978×6=5868"

Replace-Text "525×5=2625" "771×7=5397"
Replace-Text "921×7=6447" "201×7=1407"
Replace-Text "857×6=5142" "209×8=1672"
Replace-Text "611×4=2444" "675×8=5400"
Replace-Text "456×7=3192" "861×2=1722"

Replace-Text "820×5=4100" "131×4=524"
Replace-Text "433×9=3897" "403×2=806"
Replace-Text "893×6=5358" "589×8=4712"
Replace-Text "402×8=3216" "942×6=5652"
Replace-Text "881×9=7929" "901×6=5406"
